$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.099.68'
$ws.Range("E2").Value = '  +4.36%  '

# Row 3
$ws.Range("D3").Value = '3.343.81'
$ws.Range("E3").Value = '  +2.82%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '413.37'
$ws.Range("E5").Value = '  +3.78%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.47'
$ws.Range("E6").Value = '  +0.08%  '

# Row 7
$ws.Range("E7").Value = '  +4.91%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.635'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.12'
$ws.Range("E10").Value = '  +1.93%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0985'
$ws.Range("E11").Value = '  +3.93%  '

# Row 12
$ws.Range("E12").Value = '  +1.41%  '

# Row 13
$ws.Range("D13").Value = '3.884.94'
$ws.Range("E13").Value = '  +3.20%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.45'
$ws.Range("E14").Value = '  +3.92%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.48'
$ws.Range("E15").Value = '  +1.45%  '

# Row 16
$ws.Range("D16").Value = '3.330.96'
$ws.Range("E16").Value = '  +2.55%  '

# Row 17
$ws.Range("E17").Value = '  -0.12%  '

# Row 18
$ws.Range("D18").Value = '58.995.30'
$ws.Range("E18").Value = '  +4.53%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.93'
$ws.Range("E19").Value = '  -1.61%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.36'
$ws.Range("E20").Value = '  +1.29%  '

# Row 21
$ws.Range("E21").Value = '  +5.91%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.07'
$ws.Range("E22").Value = '  +0.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '303.90'
$ws.Range("E23").Value = '  +1.26%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.35'
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$ws.Range("E25").Value = '  +0.15%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.70'
$ws.Range("E26").Value = '  +1.26%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.48'
$ws.Range("E27").Value = '  +2.84%  '

# Row 28
$ws.Range("E28").Value = '  -1.65%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.45'
$ws.Range("E29").Value = '  +1.10%  '

# Row 30
$ws.Range("E30").Value = '  +0.52%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.115'
$ws.Range("E31").Value = '  +3.42%  '

# Row 32
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.63'
$ws.Range("E32").Value = '  +4.66%  '

# Row 33
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.03%  '

# Row 34
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.11'
$ws.Range("E34").Value = '  +9.33%  '

# Row 35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0537'
$ws.Range("E35").Value = '  +10.77%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.14'
$ws.Range("E36").Value = '  +0.66%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.12'
$ws.Range("E37").Value = '  +0.89%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.16'
$ws.Range("E38").Value = '  +0.97%  '

# Row 39
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.15%  '

# Row 40
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.48'
$ws.Range("E40").Value = '  -1.79%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '137.87'
$ws.Range("E41").Value = '  +2.53%  '

# Row 42
$ws.Range("E42").Value = '  +1.79%  '

# Row 43
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.04'
$ws.Range("E43").Value = '  +1.58%  '

# Row 44
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.92'
$ws.Range("E44").Value = '  -0.68%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.10'
$ws.Range("E45").Value = '  -2.76%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.281'
$ws.Range("E46").Value = '  -1.33%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.25'
$ws.Range("E47").Value = '  +7.82%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.48'
$ws.Range("E48").Value = '  +0.96%  '

# Row 49
$ws.Range("D49").Value = '2.203.41'
$ws.Range("E49").Value = '  +2.72%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.39'
$ws.Range("E50").Value = '  -1.03%  '

# Row 51
$ws.Range("E51").Value = '  -11.43%  '
